# Apply scraped-schedule update for Línea 141 (run 2026-01-16 08:16:48)
$wb = $excel.ActiveWorkbook

# ---- Sheet "LP1912": 33 row(s) changed ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:16:48"
$ws.Cells.Item(3,1).Value = "Total filas: 115"
$ws.Cells.Item(48,1).Value = "06:02:16"
$ws.Cells.Item(48,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(48,4).Value = 63
$ws.Cells.Item(49,1).Value = "05:18:23"
$ws.Cells.Item(49,3).Value = "15_ABASTO"
$ws.Cells.Item(49,4).Value = 107
$ws.Cells.Item(92,1).Value = "08:16:48"
$ws.Cells.Item(92,2).Value = "08:55"
$ws.Cells.Item(92,3).Value = "10_OLMOS"
$ws.Cells.Item(92,4).Value = 39
$ws.Cells.Item(93,1).Value = "07:14:27"
$ws.Cells.Item(93,2).Value = "09:01"
$ws.Cells.Item(93,4).Value = 107
$ws.Cells.Item(94,1).Value = "07:44:08"
$ws.Cells.Item(94,2).Value = "09:02"
$ws.Cells.Item(94,3).Value = "215A_EL PATO"
$ws.Cells.Item(94,4).Value = 78
$ws.Cells.Item(95,1).Value = "07:57:27"
$ws.Cells.Item(95,2).Value = "09:03"
$ws.Cells.Item(95,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(95,4).Value = 66
$ws.Cells.Item(96,1).Value = "08:16:48"
$ws.Cells.Item(96,2).Value = "09:04"
$ws.Cells.Item(96,4).Value = 48
$ws.Cells.Item(97,1).Value = "07:14:27"
$ws.Cells.Item(97,2).Value = "09:07"
$ws.Cells.Item(97,4).Value = 113
$ws.Cells.Item(98,1).Value = "07:44:08"
$ws.Cells.Item(98,2).Value = "09:08"
$ws.Cells.Item(98,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(98,4).Value = 84
$ws.Cells.Item(99,1).Value = "07:57:27"
$ws.Cells.Item(99,2).Value = "09:09"
$ws.Cells.Item(99,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(99,4).Value = 72
$ws.Cells.Item(100,1).Value = "07:14:27"
$ws.Cells.Item(100,2).Value = "09:10"
$ws.Cells.Item(100,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(100,4).Value = 116
$ws.Cells.Item(101,1).Value = "07:44:08"
$ws.Cells.Item(101,2).Value = "09:11"
$ws.Cells.Item(101,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(101,4).Value = 87
$ws.Cells.Item(102,2).Value = "09:14"
$ws.Cells.Item(102,3).Value = "16_SANTA ANA"
$ws.Cells.Item(102,4).Value = 90
$ws.Cells.Item(103,1).Value = "07:57:27"
$ws.Cells.Item(103,2).Value = "09:16"
$ws.Cells.Item(103,3).Value = "27_EL RETIRO"
$ws.Cells.Item(103,4).Value = 79
$ws.Cells.Item(104,1).Value = "07:44:08"
$ws.Cells.Item(104,2).Value = "09:17"
$ws.Cells.Item(104,3).Value = "27_EL RETIRO"
$ws.Cells.Item(104,4).Value = 93
$ws.Cells.Item(105,1).Value = "07:44:08"
$ws.Cells.Item(105,2).Value = "09:21"
$ws.Cells.Item(105,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(105,4).Value = 97
$ws.Cells.Item(106,2).Value = "09:22"
$ws.Cells.Item(106,3).Value = "16_SANTA ANA"
$ws.Cells.Item(106,4).Value = 85
$ws.Cells.Item(107,1).Value = "07:57:27"
$ws.Cells.Item(107,2).Value = "09:22"
$ws.Cells.Item(107,4).Value = 85
$ws.Cells.Item(108,2).Value = "09:23"
$ws.Cells.Item(108,3).Value = "17_ROMERO"
$ws.Cells.Item(108,4).Value = 99
$ws.Cells.Item(109,1).Value = "07:57:27"
$ws.Cells.Item(109,2).Value = "09:23"
$ws.Cells.Item(109,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(109,4).Value = 86
$ws.Cells.Item(110,2).Value = "09:24"
$ws.Cells.Item(110,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(110,4).Value = 100
$ws.Cells.Item(111,1).Value = "08:16:48"
$ws.Cells.Item(111,2).Value = "09:29"
$ws.Cells.Item(111,3).Value = "16_SANTA ANA"
$ws.Cells.Item(111,4).Value = 73
$ws.Cells.Item(112,2).Value = "09:32"
$ws.Cells.Item(112,3).Value = "15_ABASTO"
$ws.Cells.Item(112,4).Value = 108
$ws.Cells.Item(113,1).Value = "07:44:08"
$ws.Cells.Item(113,2).Value = "09:33"
$ws.Cells.Item(113,3).Value = "10_OLMOS"
$ws.Cells.Item(113,4).Value = 109
$ws.Cells.Item(114,1).Value = "07:44:08"
$ws.Cells.Item(114,2).Value = "09:36"
$ws.Cells.Item(114,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(114,4).Value = 112
$ws.Cells.Item(114,5).Value = "LP1912"
$ws.Cells.Item(115,1).Value = "08:16:48"
$ws.Cells.Item(115,2).Value = "09:37"
$ws.Cells.Item(115,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(115,4).Value = 81
$ws.Cells.Item(115,5).Value = "LP1912"
$ws.Cells.Item(116,1).Value = "08:16:48"
$ws.Cells.Item(116,2).Value = "09:41"
$ws.Cells.Item(116,3).Value = "215C_EL PATO"
$ws.Cells.Item(116,4).Value = 85
$ws.Cells.Item(116,5).Value = "LP1912"
$ws.Cells.Item(117,1).Value = "07:44:08"
$ws.Cells.Item(117,2).Value = "09:42"
$ws.Cells.Item(117,3).Value = "215C_EL PATO"
$ws.Cells.Item(117,4).Value = 118
$ws.Cells.Item(117,5).Value = "LP1912"
$ws.Cells.Item(118,1).Value = "07:57:27"
$ws.Cells.Item(118,2).Value = "09:43"
$ws.Cells.Item(118,3).Value = "14_ABASTO"
$ws.Cells.Item(118,4).Value = 106
$ws.Cells.Item(118,5).Value = "LP1912"
$ws.Cells.Item(119,1).Value = "08:16:48"
$ws.Cells.Item(119,2).Value = "10:10"
$ws.Cells.Item(119,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(119,4).Value = 114
$ws.Cells.Item(119,5).Value = "LP1912"
$ws.Cells.Item(120,1).Value = "08:16:48"
$ws.Cells.Item(120,2).Value = "10:12"
$ws.Cells.Item(120,3).Value = "15_ABASTO"
$ws.Cells.Item(120,4).Value = 116
$ws.Cells.Item(120,5).Value = "LP1912"

# ---- Sheet "LP1912-215": 4 row(s) changed ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:16:48"
$ws.Cells.Item(3,1).Value = "Total filas: 20"
$ws.Cells.Item(22,1).Value = "08:16:48"
$ws.Cells.Item(22,2).Value = "09:41"
$ws.Cells.Item(22,4).Value = 85
$ws.Cells.Item(23,1).Value = "07:44:08"
$ws.Cells.Item(23,2).Value = "09:42"
$ws.Cells.Item(23,3).Value = "215C_EL PATO"
$ws.Cells.Item(23,4).Value = 118
$ws.Cells.Item(23,5).Value = "LP1912"

# ---- Sheet "6203-6173": 5 row(s) changed ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:16:48"
$ws.Cells.Item(3,1).Value = "Total filas: 20"
$ws.Cells.Item(23,1).Value = "08:16:48"
$ws.Cells.Item(23,2).Value = "09:08"
$ws.Cells.Item(23,4).Value = 52
$ws.Cells.Item(24,1).Value = "07:14:27"
$ws.Cells.Item(24,2).Value = "09:09"
$ws.Cells.Item(24,3).Value = "215D_LA PLATA"
$ws.Cells.Item(24,4).Value = 115
$ws.Cells.Item(24,5).Value = "L6203"
$ws.Cells.Item(25,1).Value = "08:16:48"
$ws.Cells.Item(25,2).Value = "10:02"
$ws.Cells.Item(25,3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(25,4).Value = 106
$ws.Cells.Item(25,5).Value = "L6173"

